# Update column G ("K" = strikeouts) for rows 2-64 of Sheet1, per the
# regenerated save_data (K replaces Strike#, std/mean recomputed, s_vals
# written elsewhere). Only the G-column values change in this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 2
    9  = 0
    10 = 0
    11 = 4
    12 = 1
    13 = 2
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    19 = 3
    20 = 0
    21 = 3
    22 = 0
    23 = 0
    24 = 2
    25 = 2
    26 = 2
    27 = 2
    28 = 2
    29 = 1
    30 = 2
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 2
    36 = 1
    37 = 4
    38 = 2
    39 = 3
    40 = 1
    41 = 3
    42 = 0
    43 = 5
    44 = 4
    45 = 1
    46 = 3
    47 = 2
    48 = 3
    49 = 3
    50 = 2
    51 = 3
    52 = 2
    53 = 2
    54 = 3
    55 = 4
    56 = 1
    57 = 3
    58 = 2
    59 = 0
    60 = 2
    61 = 1
    62 = 1
    63 = 2
    64 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
